$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (Point=2)
$ws.Range("D3").Value = 1.5829186211763999
$ws.Range("E3").Value = 2.3834497378738799
$ws.Range("F3").Value = 4.8237095078883101
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0

# Row 4 (Point=3)
$ws.Range("D4").Value = 3.2328109165183601
$ws.Range("E4").Value = 4.9140901508588897
$ws.Range("F4").Value = 10.254590020082301
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0

# Row 5 (Point=4)
$ws.Range("D5").Value = 4.9578291966636696
$ws.Range("E5").Value = 7.6145033925869701
$ws.Range("F5").Value = 16.472686821167098
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 0

# Row 6 (Point=5)
$ws.Range("D6").Value = 6.7685826211661597
$ws.Range("E6").Value = 10.5141862150879
$ws.Range("F6").Value = 23.7538010996631
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 0

# Row 7 (Point=6)
$ws.Range("D7").Value = 8.6796623096080694
$ws.Range("E7").Value = 13.6534662591543
$ws.Range("F7").Value = 32.5530713763591
$ws.Range("G7").Value = 0
$ws.Range("H7").Value = 0

# Row 8 (Point=7)
$ws.Range("D8").Value = 10.712795323462499
$ws.Range("E8").Value = 17.091277005555099
$ws.Range("F8").Value = 43.706348329292098
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = 0

# Row 9 (Point=8)
$ws.Range("D9").Value = 12.9044437163115
$ws.Range("E9").Value = 20.923177180886999
$ws.Range("F9").Value = 59.029551977784898
$ws.Range("G9").Value = 0
$ws.Range("H9").Value = 0

# Row 10 (Point=9)
$ws.Range("D10").Value = 15.3289157201945
$ws.Range("E10").Value = 25.334358792755999
$ws.Range("F10").Value = 83.910268389632904
$ws.Range("G10").Value = 0
$ws.Range("H10").Value = 0

# Row 11 (Point=10)
$ws.Range("D11").Value = 18.206214778922899
$ws.Range("E11").Value = 30.841186852952099
$ws.Range("F11").Value = 160.82223283354401
$ws.Range("G11").Value = 0
$ws.Range("H11").Value = 0

# Reset the sheet's selection back to the default top-left cell
# (the source workbook had a stray selection at E15).
$ws.Range("A1").Select() | Out-Null
